$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D (Price) to text format so numeric-looking
# strings like "8.90" / "21.30" / "0.999" keep their exact text instead of
# being coerced into floating point numbers (matches the inlineStr text cells
# in the source workbook).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.396.49'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '3.163.67'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '570.98'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Value = '163.85'
$ws.Range("E6").Value = '  -3.22%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  -4.28%  '
$ws.Range("D9").Value = '0.117'
$ws.Range("E9").Value = '  -3.20%  '
$ws.Range("D10").Value = '6.63'
$ws.Range("E10").Value = '  -1.38%  '
$ws.Range("D11").Value = '0.385'
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").Value = '3.715.37'
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").Value = '0.127'
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").Value = '64.468.54'
$ws.Range("D15").Value = '25.31'
$ws.Range("E15").Value = '  -0.59%  '
$ws.Range("D16").Value = '3.162.76'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("E17").Value = '  -2.48%  '
$ws.Range("D18").Value = '409.27'
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("D19").Value = '12.73'
$ws.Range("E19").Value = '  -0.98%  '
$ws.Range("D20").Value = '5.27'
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("D21").Value = '7.09'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").Value = '68.67'
$ws.Range("E23").Value = '  -2.87%  '
$ws.Range("B24").Value = 'Kaspa'
$ws.Range("C24").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D24").Value = '0.198'
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").Value = '0.483'
$ws.Range("E25").Value = '  -2.27%  '
$ws.Range("E26").Value = '  -6.55%  '
$ws.Range("D27").Value = '8.90'
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +1.12%  '
$ws.Range("E29").Value = '  -2.29%  '
$ws.Range("D30").Value = '21.23'
$ws.Range("E30").Value = '  -3.45%  '
$ws.Range("D31").Value = '4.91'
$ws.Range("E31").Value = '  -2.37%  '
$ws.Range("D32").Value = '6.34'
$ws.Range("E32").Value = '  -1.08%  '
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("D34").Value = '155.86'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("E35").Value = '  -2.41%  '
$ws.Range("D36").Value = '2.688.78'
$ws.Range("E36").Value = '  -2.08%  '
$ws.Range("D37").Value = '1.69'
$ws.Range("E37").Value = '  -0.79%  '
$ws.Range("D38").Value = '24.08'
$ws.Range("E38").Value = '  -4.22%  '
$ws.Range("D39").Value = '4.09'
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("E40").Value = '  -2.94%  '
$ws.Range("D41").Value = '0.0620'
$ws.Range("E41").Value = '  -1.36%  '
$ws.Range("E42").Value = '  -4.66%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.0258'
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("D44").Value = '21.53'
$ws.Range("E44").Value = '  -2.67%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = '291.69'
$ws.Range("E45").Value = '  -2.19%  '
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("D47").Value = '0.0985'
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("E48").Value = '  -7.77%  '
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").Value = '5.72'
$ws.Range("E50").Value = '  -1.80%  '
$ws.Range("D51").Value = '0.881'
$ws.Range("E51").Value = '  -5.81%  '

# Remove the temporary text formatting again so no stray style/number-format
# is left attached to the cells (restores original default styling).
$ws.Range("D2:D51").ClearFormats()
